# Loan RBI, Variable Instalments
# Insert a new (blank) column before column N on the "Repayment Schedule"
# sheet. This shifts the old "Late" column (N) one position right to O,
# and the old "Outstanding" column (P) one position right to Q, leaving a
# new blank column N in between - matching the target XML diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make sure we are editing the "Repayment Schedule" sheet (the active tab).
if ($ws.Name -ne "Repayment Schedule") {
    $ws = $wb.Worksheets.Item("Repayment Schedule")
}

# Insert a new blank column at N - everything from N onward (N..P) shifts
# one column to the right (N->O, O->P, P->Q), carrying values/styles along.
$ws.Columns("N").Insert()

# The newly inserted column picks up the width of its left neighbour (M),
# matching the width the new column received in the authored workbook.
$ws.Columns("N").ColumnWidth = $ws.Columns("M").ColumnWidth()

# Restore the selection left behind after the edit.
$ws.Range("S9").Select() | Out-Null
